$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 405, pushing the existing rows 405:425 down to 406:426
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the new weekly price-report record
$ws.Cells.Item(405, 1).Value = 4
$ws.Cells.Item(405, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(405, 3).Value = "Los Lagos"
$ws.Cells.Item(405, 4).Value = 45008
$ws.Cells.Item(405, 5).Value = 10
$ws.Cells.Item(405, 6).Value = 100112037
$ws.Cells.Item(405, 7).Value = "Cebollín"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value = 70
$ws.Cells.Item(405, 11).Value = 6500
$ws.Cells.Item(405, 12).Value = 7500
$ws.Cells.Item(405, 13).Value = 7000
$ws.Cells.Item(405, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(405, 15).Value = "Región Metropolitana"
$ws.Cells.Item(405, 16).Value = 194
$ws.Cells.Item(405, 17).Value = 36
$ws.Cells.Item(405, 18).Value = "Hortaliza"
